# DONE NEST AUTOMATION(VER 1)
# Update HIGH/LOW/CLOSE/LTP/VOL/9:25 CLOSE data for each symbol row on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns: Row, B(HIGH), C(LOW), D(CLOSE), E(LTP), F(VOL), G(9:25 CLOSE)
$data = @(
    @(2, 828.75, 806.75, 821.4, 821.95, 142, 807.9),
    @(3, 839.9, 828, 829, 830.25, 53, 833.05),
    @(4, 44626, 44443.15, 44516.4, 44500.5, 18, 44607),
    @(5, 328.65, 325.5, 327.2, 327, 143, 328.65),
    @(6, 498.35, 488.45, 494.5, 493.95, 230, 489.1),
    @(7, 459.15, 448.5, 457.35, 456.55, 108, 450.1),
    @(8, 975.2, 967, 968.1, 968.45, 273, 974.15),
    @(9, 679.9, 657.75, 675, 675.35, 90, 660.6),
    @(10, 19366.95, 19306, 19345.75, 19337.1, 49, 19344.9),
    @(11, 2442, 2410.5, 2423.15, 2422.45, 110, 2433.05),
    @(12, 575.8, 572, 575, 574.8, 318, 573.15),
    @(13, 841.5, 834.75, 836.8, 836.6, 29, 838.9),
    @(14, 611.25, 605.2, 608.25, 607.65, 172, 606.6),
    @(15, 120.35, 117.8, 120, 119.7, 950, 118.2),
    @(16, 3384.7, 3362.25, 3370.15, 3370.55, 28, 3379.2),
    @(17, 3086.45, 3046.05, 3075.4, 3074.7, 20, 3050.55)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $ws.Cells.Item($row, 7).Value = $entry[6]
}
